$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 12011.643
$ws.Range("I9").Value = 15113.333
$ws.Range("J9").Value = 9685.375
$ws.Range("K9").Value = 15113.333
$ws.Range("L9").Value = 9685.375
$ws.Range("M9").Value = -14944.333
$ws.Range("N9").Value = -10023.375
$ws.Range("H12").Value = 209.36363
$ws.Range("I12").Value = 89.71429000000001
$ws.Range("K12").Value = 89.71429000000001
$ws.Range("M12").Value = 80.28570999999999
$ws.Range("H28").Value = 2086.08
$ws.Range("I28").Value = 2272.8823
$ws.Range("K28").Value = 2272.8823
$ws.Range("M28").Value = -1787.8823
$ws.Range("H40").Value = 1344.3
$ws.Range("I40").Value = 998
$ws.Range("J40").Value = 1575.1666
$ws.Range("K40").Value = 998
$ws.Range("L40").Value = 1575.1666
$ws.Range("M40").Value = -823
$ws.Range("N40").Value = -1925.1666
$ws.Range("H43").Value = 1104.5
$ws.Range("J43").Value = 1104.5
$ws.Range("L43").Value = 1104.5
$ws.Range("N43").Value = -1242.5
$ws.Range("H80").Value = 1724.8
$ws.Range("I80").Value = 1362.25
$ws.Range("J80").Value = 3175
$ws.Range("K80").Value = 4086.75
$ws.Range("L80").Value = 9525
$ws.Range("M80").Value = -3088.75
$ws.Range("N80").Value = -11521
$ws.Range("H83").Value = 1724.8
$ws.Range("I83").Value = 1362.25
$ws.Range("J83").Value = 3175
$ws.Range("K83").Value = 12260.25
$ws.Range("L83").Value = 28575
$ws.Range("M83").Value = -7268.25
$ws.Range("N83").Value = -38559
$ws.Range("H88").Value = 1365.375
$ws.Range("I88").Value = 891.3333
$ws.Range("J88").Value = 1649.8
$ws.Range("K88").Value = 891.3333
$ws.Range("L88").Value = 1649.8
$ws.Range("M88").Value = -485.3333
$ws.Range("N88").Value = -2461.8
$ws.Range("H91").Value = 1365.375
$ws.Range("I91").Value = 891.3333
$ws.Range("J91").Value = 1649.8
$ws.Range("K91").Value = 891.3333
$ws.Range("L91").Value = 1649.8
$ws.Range("M91").Value = 512.6667
$ws.Range("N91").Value = -4457.8
$ws.Range("H92").Value = 10005957
$ws.Range("I92").Value = 14293279
$ws.Range("K92").Value = 14293279
$ws.Range("M92").Value = -14292031
$ws.Range("H135").Value = 18524804
$ws.Range("I135").Value = 19614218
$ws.Range("J135").Value = 4751
$ws.Range("K135").Value = 176527962
$ws.Range("L135").Value = 42759
$ws.Range("M135").Value = -176525427
$ws.Range("N135").Value = -47829

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2744.5518
$ws.Range("I63").Value = 2744.5518
$ws.Range("K63").Value = 2744.5518
$ws.Range("M63").Value = -2058.5518
$ws.Range("H66").Value = 2744.5518
$ws.Range("I66").Value = 2744.5518
$ws.Range("K66").Value = 13722.759
$ws.Range("M66").Value = -10290.759
$ws.Range("H74").Value = 17248034
$ws.Range("I74").Value = 3519.75
$ws.Range("K74").Value = 3519.75
$ws.Range("M74").Value = -2645.75
$ws.Range("H77").Value = 17248034
$ws.Range("I77").Value = 3519.75
$ws.Range("K77").Value = 17598.75
$ws.Range("M77").Value = -13230.75
$ws.Range("H88").Value = 3214
$ws.Range("I88").Value = 3078.2222
$ws.Range("K88").Value = 3078.2222
$ws.Range("M88").Value = -2672.2222
$ws.Range("H91").Value = 3214
$ws.Range("I91").Value = 3078.2222
$ws.Range("K91").Value = 3078.2222
$ws.Range("M91").Value = -1674.2222
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -59178
$ws.Range("H132").Value = 2785.0322
$ws.Range("I132").Value = 2129.2083
$ws.Range("K132").Value = 6387.624899999999
$ws.Range("M132").Value = -3857.624899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 214.85715
$ws.Range("I7").Value = 78.55556
$ws.Range("K7").Value = 78.55556
$ws.Range("M7").Value = 34.44444
$ws.Range("H8").Value = 5000
$ws.Range("J8").Value = 5000
$ws.Range("L8").Value = 5000
$ws.Range("N8").Value = -5280
$ws.Range("J59").Value = 26000
$ws.Range("L59").Value = 26000
$ws.Range("N59").Value = -28290
$ws.Range("H132").Value = 1294
$ws.Range("I132").Value = 1294
$ws.Range("K132").Value = 3882
$ws.Range("M132").Value = -1352
$ws.Range("H134").Value = 11589.472
$ws.Range("I134").Value = 11984.9
$ws.Range("K134").Value = 35954.7
$ws.Range("M134").Value = -33419.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 838.6923
$ws.Range("I129").Value = 440.3
$ws.Range("K129").Value = 1320.9
$ws.Range("M129").Value = 3679.1
$ws.Range("H131").Value = 23635.191
$ws.Range("J131").Value = 2347.2104
$ws.Range("L131").Value = 7041.6312
$ws.Range("N131").Value = -17121.6312

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 27812276
$ws.Range("I18").Value = 37049700
$ws.Range("K18").Value = 37049700
$ws.Range("M18").Value = -37049407
$ws.Range("H123").Value = 25914.285
$ws.Range("J123").Value = 25900
$ws.Range("L123").Value = 25900
$ws.Range("N123").Value = -30800
$ws.Range("H132").Value = 44369
$ws.Range("I132").Value = 44369
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 133107
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -130577
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3240.5881
$ws.Range("I16").Value = 3526.7273
$ws.Range("J16").Value = 2716
$ws.Range("K16").Value = 3526.7273
$ws.Range("L16").Value = 2716
$ws.Range("M16").Value = -3356.7273
$ws.Range("N16").Value = -3056
$ws.Range("H22").Value = 1840.8667
$ws.Range("I22").Value = 1395.8
$ws.Range("J22").Value = 2063.4
$ws.Range("K22").Value = 1395.8
$ws.Range("L22").Value = 2063.4
$ws.Range("M22").Value = -1100.8
$ws.Range("N22").Value = -2653.4
$ws.Range("H27").Value = 1840.8667
$ws.Range("I27").Value = 1395.8
$ws.Range("J27").Value = 2063.4
$ws.Range("K27").Value = 1395.8
$ws.Range("L27").Value = 2063.4
$ws.Range("M27").Value = -1288.8
$ws.Range("N27").Value = -2277.4
$ws.Range("H46").Value = 1554
$ws.Range("I46").Value = 1554
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1554
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1366
$ws.Range("N46").ClearContents()
$ws.Range("H61").Value = 5256.1904
$ws.Range("I61").Value = 4828.8237
$ws.Range("J61").Value = 7072.5
$ws.Range("K61").Value = 4828.8237
$ws.Range("L61").Value = 7072.5
$ws.Range("M61").Value = -4626.8237
$ws.Range("N61").Value = -7476.5
$ws.Range("H68").Value = 2151.9
$ws.Range("I68").Value = 2314.0588
$ws.Range("J68").Value = 1233
$ws.Range("K68").Value = 2314.0588
$ws.Range("L68").Value = 1233
$ws.Range("M68").Value = -1565.0588
$ws.Range("N68").Value = -2731
$ws.Range("H71").Value = 2151.9
$ws.Range("I71").Value = 2314.0588
$ws.Range("J71").Value = 1233
$ws.Range("K71").Value = 11570.294
$ws.Range("L71").Value = 6165
$ws.Range("M71").Value = -7826.293999999998
$ws.Range("N71").Value = -13653
$ws.Range("H82").Value = 3765.96
$ws.Range("I82").Value = 3143.6924
$ws.Range("J82").Value = 4440.0835
$ws.Range("K82").Value = 3143.6924
$ws.Range("L82").Value = 4440.0835
$ws.Range("M82").Value = -2782.6924
$ws.Range("N82").Value = -5162.0835
$ws.Range("H85").Value = 3765.96
$ws.Range("I85").Value = 3143.6924
$ws.Range("J85").Value = 4440.0835
$ws.Range("K85").Value = 3143.6924
$ws.Range("L85").Value = 4440.0835
$ws.Range("M85").Value = -1895.6924
$ws.Range("N85").Value = -6936.0835
$ws.Range("H113").Value = 5256.1904
$ws.Range("I113").Value = 4828.8237
$ws.Range("J113").Value = 7072.5
$ws.Range("K113").Value = 4828.8237
$ws.Range("L113").Value = 7072.5
$ws.Range("M113").Value = -2658.8237
$ws.Range("N113").Value = -11412.5
$ws.Range("H136").Value = 4578
$ws.Range("I136").Value = 3854.0833
$ws.Range("K136").Value = 11562.2499
$ws.Range("M136").Value = -9012.249899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1661.8182
$ws.Range("I136").Value = 1433
$ws.Range("J136").Value = 3950
$ws.Range("K136").Value = 4299
$ws.Range("L136").Value = 11850
$ws.Range("M136").Value = -1749
$ws.Range("N136").Value = -16950
